$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "320018179991"
$ws.Range("C2").Style = "Normal"

$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "320018180002"
$ws.Range("C3").Style = "Normal"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "320018180035"
$ws.Range("C4").Style = "Normal"

$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "320018180057"
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320018180057"
$ws.Range("D5").Style = "Normal"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "320018180090"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "320018180090"
$ws.Range("D6").Style = "Normal"

$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "320018180127"
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "320018180127"
$ws.Range("D7").Style = "Normal"

$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "320018180150"
$ws.Range("C8").Style = "Normal"

$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "320018180171"
$ws.Range("C9").Style = "Normal"

$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "320018180208"
$ws.Range("C10").Style = "Normal"

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "320018180220"
$ws.Range("C11").Style = "Normal"

$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "320018180263"
$ws.Range("C12").Style = "Normal"

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "320018180285"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "320018180285"
$ws.Range("D13").Style = "Normal"

$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "320018180311"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "320018180311"
$ws.Range("D14").Style = "Normal"

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "320018180333"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "320018180333"
$ws.Range("D15").Style = "Normal"

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "320018180366"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "320018180366"
$ws.Range("D16").Style = "Normal"

$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "320018180388"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "320018180388"
$ws.Range("D17").Style = "Normal"

$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "320018180425"
$ws.Range("C18").Style = "Normal"

$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "320018180447"
$ws.Range("C19").Style = "Normal"

$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "320018180480"
$ws.Range("C20").Style = "Normal"

$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "320018180506"
$ws.Range("C21").Style = "Normal"

$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "320018180539"
$ws.Range("C22").Style = "Normal"
